$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old row 3 (second data row) so only the header + one data row remain
$ws.Rows.Item(3).Delete()

# Drop the old column U (Crop Recommendations) - new layout only spans A:T
$ws.Columns.Item(21).Delete()

# Row 1
$ws.Range("A1").Value = "Test ID"
$ws.Range("B1").Value = "Collection Date"
$ws.Range("C1").Value = "Latitude"
$ws.Range("D1").Value = "Longitude"
$ws.Range("E1").Value = "Name"
$ws.Range("F1").Value = "Area (ha)"
$ws.Range("G1").Value = "Gender"
$ws.Range("H1").Value = "Age"
$ws.Range("I1").Value = "Address"
$ws.Range("J1").Value = "Mobile No."
$ws.Range("K1").Value = "Soil pH"
$ws.Range("L1").Value = "Nitrogen"
$ws.Range("M1").Value = "Phosphorus"
$ws.Range("N1").Value = "Potassium"
$ws.Range("O1").Value = "Electrical Conductivity"
$ws.Range("P1").Value = "Temperature"
$ws.Range("Q1").Value = "Moisture"
$ws.Range("R1").Value = "Humidity"
$ws.Range("S1").Value = "Soil Health Score"
$ws.Range("T1").Value = "Recommendations"

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2222"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "12-03-2024"
$ws.Range("C2").Value = 253
$ws.Range("D2").Value = 256
$ws.Range("E2").Value = "asdasdasdasdasd"
$ws.Range("F2").Value = 25
$ws.Range("G2").Value = "Male"
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = "26a5sd5646as5d"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2323265589"
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 200
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 30
$ws.Range("R2").Value = 22
$ws.Range("S2").Value = 0.5473459137758564
$ws.Range("T2").Value = "Millets(Pearl Millet, Sorghum), Maize, Soybean, Groundnut"
